$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.928.13'
$ws.Cells.Item(2, 5).Value = '  -1.09%  '
$ws.Cells.Item(3, 4).Value = '2.192.05'
$ws.Cells.Item(3, 5).Value = '  -2.37%  '
$ws.Cells.Item(4, 5).Value = '  -0.17%  '
$ws.Cells.Item(5, 4).Value = "'294.48"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -4.44%  '
$ws.Cells.Item(6, 4).Value = "'89.04"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -5.88%  '
$ws.Cells.Item(7, 4).Value = "'0.564"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -1.22%  '
$ws.Cells.Item(8, 5).Value = '  -0.06%  '
$ws.Cells.Item(9, 4).Value = "'0.483"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -8.21%  '
$ws.Cells.Item(10, 4).Value = "'32.03"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -8.33%  '
$ws.Cells.Item(11, 4).Value = "'0.0772"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -4.85%  '
$ws.Cells.Item(12, 5).Value = '  -1.27%  '
$ws.Cells.Item(13, 5).Value = '  -5.53%  '
$ws.Cells.Item(14, 4).Value = '2.525.86'
$ws.Cells.Item(14, 5).Value = '  -2.47%  '
$ws.Cells.Item(15, 4).Value = '2.256.29'
$ws.Cells.Item(15, 5).Value = '  -4.18%  '
$ws.Cells.Item(16, 4).Value = "'13.13"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -3.96%  '
$ws.Cells.Item(17, 5).Value = '  -8.02%  '
$ws.Cells.Item(18, 4).Value = '43.597.55'
$ws.Cells.Item(18, 5).Value = '  -1.16%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0889'
$ws.Cells.Item(19, 5).Value = '  -7.95%  '
$ws.Cells.Item(20, 5).Value = '  -8.93%  '
$ws.Cells.Item(21, 4).Value = "'10.81"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -12.04%  '
$ws.Cells.Item(22, 4).Value = "'63.10"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -4.16%  '
$ws.Cells.Item(23, 4).Value = "'231.65"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -2.38%  '
$ws.Cells.Item(24, 5).Value = '  -9.32%  '
$ws.Cells.Item(25, 5).Value = '  +0.74%  '
$ws.Cells.Item(26, 5).Value = '  -8.94%  '
$ws.Cells.Item(27, 5).Value = '  +0.55%  '
$ws.Cells.Item(28, 4).Value = "'36.27"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -5.57%  '
$ws.Cells.Item(29, 4).Value = "'9.24"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -6.28%  '
$ws.Cells.Item(30, 4).Value = "'19.21"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -4.28%  '
$ws.Cells.Item(31, 4).Value = "'148.11"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -3.88%  '
$ws.Cells.Item(32, 5).Value = '  -11.59%  '
$ws.Cells.Item(33, 5).Value = '  -5.22%  '
$ws.Cells.Item(34, 4).Value = "'0.0734"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -8.13%  '
$ws.Cells.Item(35, 5).Value = '  -3.98%  '
$ws.Cells.Item(36, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(36, 4).Value = "'2.82"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -9.05%  '
$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(37, 4).Value = "'0.102"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -6.04%  '
$ws.Cells.Item(38, 4).Value = "'1.64"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -8.97%  '
$ws.Cells.Item(39, 4).Value = "'0.0284"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -6.88%  '
$ws.Cells.Item(40, 4).Value = "'3.51"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -8.25%  '
$ws.Cells.Item(41, 5).Value = '  -11.83%  '
$ws.Cells.Item(42, 4).Value = "'1.00"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -0.34%  '
$ws.Cells.Item(43, 4).Value = "'13.01"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -11.56%  '
$ws.Cells.Item(44, 4).Value = '1.796.18'
$ws.Cells.Item(44, 5).Value = '  +3.11%  '
$ws.Cells.Item(45, 4).Value = "'1.66"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +4.09%  '
$ws.Cells.Item(46, 5).Value = '  +11.50%  '
$ws.Cells.Item(47, 4).Value = "'0.174"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -10.37%  '
$ws.Cells.Item(48, 4).Value = "'72.28"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -10.32%  '
$ws.Cells.Item(49, 4).Value = "'91.70"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -8.14%  '
$ws.Cells.Item(50, 4).Value = "'64.66"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -8.61%  '
$ws.Cells.Item(51, 4).Value = '2.409.23'
$ws.Cells.Item(51, 5).Value = '  -2.39%  '
